$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values of row 2 and row 3 for columns D, M, N, O, P, R, S
$cols = @("D", "M", "N", "O", "P", "R", "S")

foreach ($col in $cols) {
    $addr2 = "$col" + "2"
    $addr3 = "$col" + "3"
    $val2 = $ws.Range($addr2).Value2
    $val3 = $ws.Range($addr3).Value2
    $ws.Range($addr2).Value2 = $val3
    $ws.Range($addr3).Value2 = $val2
}
